$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID Visita"
$ws.Range("B1").Value = "Establecimiento"
$ws.Range("C1").Value = "Tipo de documento"
$ws.Range("D1").Value = "Nro documento"
$ws.Range("E1").Value = "Uso del tapabocas"
$ws.Range("F1").Value = "Temperatura"
$ws.Range("G1").Value = "Fecha de ingreso"
$ws.Range("H1").Value = "Hora de ingreso"
$ws.Range("I1").Value = "Ingreso"
$ws.Range("J1").Value = "Nombres"
$ws.Range("K1").Value = "Apellidos"

# The new header cells (J1, K1) didn't exist before, so they have no style yet.
# Copy the header style from an existing header cell onto them.
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "PepsiCo, Inc."
$ws.Range("C2").Value = "Cédula de ciudadanía"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1339998889"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "Si"
$ws.Range("F2").Value = 35.7
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2020-12-05"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = "02:38:13.203773"
$ws.Range("I2").Value = "Aceptado"
$ws.Range("J2").Value = "Anuel"
$ws.Range("K2").Value = "AA"

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "PepsiCo, Inc."
$ws.Range("C3").Value = "Cédula de ciudadanía"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1193474912"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = 36
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2020-12-05"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = "02:40:18.030912"
$ws.Range("I3").Value = "Denegado"
$ws.Range("J3").Value = "Isabela"
$ws.Range("K3").Value = "Acevedo García"

$wb.Save()
